$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 updates
$ws.Range("Q2").Value = 2.08
$ws.Range("R2").Value = 1.73

# Row 5 updates
$ws.Range("G5").Value = 3.6
$ws.Range("I5").Value = 2.15
$ws.Range("J5").Value = 4
$ws.Range("AA5").Value = 29
$ws.Range("AH5").Value = 7

# Row 6 updates
$ws.Range("M6").Value = 1.13
$ws.Range("N6").Value = 6

# Row 8 updates
$ws.Range("O8").Value = 1.36
$ws.Range("P8").Value = 3
